$d = $word.ActiveDocument

# --- Step 1: delete the long tail of paragraphs (old #18 through the last one) ---
$total = $d.Paragraphs.Count
if ($total -ge 18) {
    $delStart = $d.Paragraphs.Item(18).Range.Start
    $delEnd = $d.Paragraphs.Item($total).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

# --- Step 2: rewrite the text of the remaining paragraphs (old #1 heading + #2-17 rows) ---
$newTexts = @(
    "🧾 Deposition Topic Table of Contents",
    "- **Page 1** · Page 1 · Line 1",
    "- **Q: Please state your name for the record.** · Page 1 · Line 3",
    "- **Q: Please state your name for the record.** · Page 1 · Line 4",
    "- **Q: Are you currently employed?** · Page 1 · Line 6",
    "- **A: Yes, I do.** · Page 1 · Line 7",
    "- **Q: Are you currently employed?** · Page 1 · Line 9",
    "- **A: Yes, I work at Horizon Corp as a senior analyst.** · Page 1 · Line 10",
    "- **Q: Are you currently employed?** · Page 1 · Line 12",
    "- **A: Yes, I do.** · Page 1 · Line 13",
    "- **Page 1** · Page 1 · Line 15",
    "- **A: Yes, I work at Horizon Corp as a senior analyst.** · Page 1 · Line 17",
    "- **A: Yes, I work at Horizon Corp as a senior analyst.** · Page 1 · Line 18",
    "- **Q: Are you currently employed?** · Page 1 · Line 20",
    "- **A: Yes, I do.** · Page 1 · Line 21",
    "- **Q: Please state your name for the record.** · Page 1 · Line 23",
    "- **A: Yes, I do.** · Page 1 · Line 24"
)

for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $p = $d.Paragraphs.Item($i + 1)
    $p.Range.Text = $newTexts[$i]
}

# --- Step 3: insert a blank paragraph right after the heading (before old row #2,
#     which carries no special style) so the new blank paragraph stays plain ---
$firstRow = $d.Paragraphs.Item(2)
$firstRow.Range.InsertParagraphBefore()

# --- Step 4: insert a blank paragraph at the very end of the document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
